# Scheduled market-data refresh: updates derived price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
#  LeveProfitNQ/HQ -- columns H:N) on the per-class Leve profit sheets
# using the latest pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3255.9333
$ws.Range("I51").Value = 2271.2856
$ws.Range("J51").Value = 4117.5
$ws.Range("K51").Value = 2271.2856
$ws.Range("L51").Value = 4117.5
$ws.Range("M51").Value = -1787.2856
$ws.Range("N51").Value = -5085.5

$ws.Range("H137").Value = 7730.923
$ws.Range("I137").Value = 7850.2
$ws.Range("J137").Value = 7333.3335
$ws.Range("K137").Value = 23550.6
$ws.Range("L137").Value = 22000.0005
$ws.Range("M137").Value = -21000.6
$ws.Range("N137").Value = -27100.0005

$ws.Range("H138").Value = 2299.5334
$ws.Range("I138").Value = 1281.3549
$ws.Range("J138").Value = 3387.9312
$ws.Range("K138").Value = 3844.0647
$ws.Range("L138").Value = 10163.7936
$ws.Range("M138").Value = 1295.9353
$ws.Range("N138").Value = -20443.7936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23965.021
$ws.Range("I32").Value = 5711.7637
$ws.Range("K32").Value = 5711.7637
$ws.Range("M32").Value = -5424.7637

$ws.Range("H61").Value = 1867.742
$ws.Range("I61").Value = 1495.6957
$ws.Range("J61").Value = 2937.375
$ws.Range("K61").Value = 1495.6957
$ws.Range("L61").Value = 2937.375
$ws.Range("M61").Value = -1283.6957
$ws.Range("N61").Value = -3361.375

$ws.Range("H122").Value = 1326.091
$ws.Range("I122").Value = 817.125
$ws.Range("J122").Value = 2683.3333
$ws.Range("K122").Value = 2451.375
$ws.Range("L122").Value = 8049.999899999999
$ws.Range("M122").Value = -1.375
$ws.Range("N122").Value = -12949.9999

$ws.Range("H136").Value = 1867.742
$ws.Range("I136").Value = 1495.6957
$ws.Range("J136").Value = 2937.375
$ws.Range("K136").Value = 4487.0871
$ws.Range("L136").Value = 8812.125
$ws.Range("M136").Value = -1937.0871
$ws.Range("N136").Value = -13912.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1600.625
$ws.Range("I99").Value = 1473.6364
$ws.Range("J99").Value = 1880
$ws.Range("K99").Value = 1473.6364
$ws.Range("L99").Value = 1880
$ws.Range("M99").Value = 24.36359999999991
$ws.Range("N99").Value = -4876

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1965.9524
$ws.Range("I132").Value = 1573.1666
$ws.Range("J132").Value = 4322.6665
$ws.Range("K132").Value = 4719.4998
$ws.Range("L132").Value = 12967.9995
$ws.Range("M132").Value = -2189.4998
$ws.Range("N132").Value = -18027.9995

$ws.Range("H134").Value = 2070.4285
$ws.Range("I134").Value = 1835.8422
$ws.Range("J134").Value = 2565.6667
$ws.Range("K134").Value = 5507.5266
$ws.Range("L134").Value = 7697.000100000001
$ws.Range("M134").Value = -2972.5266
$ws.Range("N134").Value = -12767.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34.52381
$ws.Range("I2").Value = 29.833334
$ws.Range("J2").Value = 40.77778
$ws.Range("K2").Value = 179.000004
$ws.Range("L2").Value = 244.66668
$ws.Range("M2").Value = -66.00000399999999
$ws.Range("N2").Value = -470.66668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1900
$ws.Range("I80").Value = 1850
$ws.Range("J80").Value = 1950
$ws.Range("K80").Value = 1850
$ws.Range("L80").Value = 1950
$ws.Range("M80").Value = -852
$ws.Range("N80").Value = -3946

$ws.Range("H83").Value = 1900
$ws.Range("I83").Value = 1850
$ws.Range("J83").Value = 1950
$ws.Range("K83").Value = 9250
$ws.Range("L83").Value = 9750
$ws.Range("M83").Value = -4258
$ws.Range("N83").Value = -19734

$ws.Range("H123").Value = 18117.666
$ws.Range("J123").Value = 18117.666
$ws.Range("L123").Value = 18117.666
$ws.Range("N123").Value = -23017.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1655.5
$ws.Range("I7").Value = 1154.8889
$ws.Range("J7").Value = 2156.111
$ws.Range("K7").Value = 1154.8889
$ws.Range("L7").Value = 2156.111
$ws.Range("M7").Value = -1042.8889
$ws.Range("N7").Value = -2380.111

$ws.Range("H22").Value = 425.5
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 451
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 451
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1041

$ws.Range("H27").Value = 425.5
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 451
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 451
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -665

$ws.Range("H55").Value = 307.66666
$ws.Range("I55").Value = 340.30768
$ws.Range("J55").Value = 95.5
$ws.Range("K55").Value = 340.30768
$ws.Range("L55").Value = 95.5
$ws.Range("M55").Value = -167.30768
$ws.Range("N55").Value = -441.5

$ws.Range("H93").Value = 940.4583
$ws.Range("I93").Value = 707.13336
$ws.Range("J93").Value = 1329.3334
$ws.Range("K93").Value = 707.13336
$ws.Range("L93").Value = 1329.3334
$ws.Range("M93").Value = 540.86664
$ws.Range("N93").Value = -3825.3334

$ws.Range("H126").Value = 1655.5
$ws.Range("I126").Value = 1154.8889
$ws.Range("J126").Value = 2156.111
$ws.Range("K126").Value = 3464.6667
$ws.Range("L126").Value = 6468.333
$ws.Range("M126").Value = -994.6666999999998
$ws.Range("N126").Value = -11408.333

$ws.Range("H132").Value = 6178.3257
$ws.Range("I132").Value = 4303.048
$ws.Range("J132").Value = 7968.364
$ws.Range("K132").Value = 12909.144
$ws.Range("L132").Value = 23905.092
$ws.Range("M132").Value = -10379.144
$ws.Range("N132").Value = -28965.092

$ws.Range("H136").Value = 3245.3333
$ws.Range("I136").Value = 2913.9312
$ws.Range("J136").Value = 3984.6155
$ws.Range("K136").Value = 8741.793600000001
$ws.Range("L136").Value = 11953.8465
$ws.Range("M136").Value = -6191.793600000001
$ws.Range("N136").Value = -17053.8465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1319.9565
$ws.Range("I132").Value = 796.7895
$ws.Range("J132").Value = 3805
$ws.Range("K132").Value = 2390.3685
$ws.Range("L132").Value = 11415
$ws.Range("M132").Value = 139.6315
$ws.Range("N132").Value = -16475

$ws.Range("H136").Value = 956.8261
$ws.Range("I136").Value = 1033.1333
$ws.Range("J136").Value = 813.75
$ws.Range("K136").Value = 3099.3999
$ws.Range("L136").Value = 2441.25
$ws.Range("M136").Value = -549.3998999999999
$ws.Range("N136").Value = -7541.25
